# Updated cryptos list on Tue Aug 22 21:59:23 UTC 2023 with GitHub Actions
# Applies updated price / 1h-volume figures (and the Maker/FraxShare row swap)
# for the cryptos.xlsx source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'25.741.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.97%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.608.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.05%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'207.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.97%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'0.5184"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.72%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.2548"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.17%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.06179"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.88%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'20.04"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -6.34%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.07518"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.67%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'1.604.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -4.30%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'4.336"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.14%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'1.835.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.72%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'0.5386"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.53%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.0₅7780"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.38%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'63.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -5.33%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'25.746.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.21%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'  +0.09%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'4.573"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -5.35%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'182.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.11%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'9.948"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -4.71%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = "'  +0.03%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'5.992"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.59%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'144.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.78%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  -4.80%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'7.294"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -4.00%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'15.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.89%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'1.354"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.63%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'0.05887"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -5.47%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'1.233"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.16%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'3.348"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.72%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'3.298"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.38%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'1.580"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.51%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'0.9575"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.75%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'2.380"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.10%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'2.695"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.41%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.5664"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -6.85%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.01577"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.95%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").Value = "'  -0.42%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.8323"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -4.57%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("B42").Value = "'FraxShare"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'5.621"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -8.03%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = "'Maker"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'1.016.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -8.34%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'98.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.19%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'1.761.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.58%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "'  -1.28%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.9992"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.64%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'53.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -4.42%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.05153"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.59%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'7.808"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.53%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.4211"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.07%  "
$ws.Range("E51").Style = "Normal"

